# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows at the top of the data block (row 10),
# pushing the existing rows 10-71 down to 15-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 10-71 down to 15-76 by inserting 5 new rows at row 10.
$ws.Rows("10:14").Insert()

# Common (unchanged across the whole data block) column values.
$mercado = "Agrícola del Norte S.A. de Arica"
$region  = "Arica y Parinacota"
$codreg  = 15
$tipo    = "Fruta"
$prodId  = 100103
$prod    = "Frutos de hueso (carozo)"
$catId   = 100103004
$cat     = "Durazno"

function Set-DataRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value  = 1
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $prodId
    $ws.Cells.Item($Row, 8).Value  = $prod
    $ws.Cells.Item($Row, 9).Value  = $catId
    $ws.Cells.Item($Row, 10).Value = $cat
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-DataRow 10 44971 "Doctor Davis"  "Primera" 150 2600  27000 10733 "$/bandeja 18 kilos granel" "Región de O'Higgins" 596  18
Set-DataRow 11 44971 "Doctor Davis"  "Segunda" 150 23000 24000 23667 "$/caja 18 kilos granel"    "Región de O'Higgins" 1315 18
Set-DataRow 12 44971 "September Sun" "Primera" 250 26000 27000 26400 "$/bandeja 18 kilos granel" "Región de O'Higgins" 1467 18
Set-DataRow 13 44971 "September Sun" "Segunda" 150 23000 24000 23667 "$/caja 18 kilos granel"    "Región de O'Higgins" 1315 18
Set-DataRow 14 44971 "September Sun" "Tercera" 170 19000 20000 19529 "$/caja 18 kilos granel"    "Región de O'Higgins" 1085 18
